# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.250.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.93%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.662.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.12%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.33%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'219.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.12%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5225"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.81%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.37%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2670"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.57%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06348"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.80%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.78%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07731"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.79%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'Polkadot"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'4.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.40%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.647.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.34%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.891.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.94%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5478"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.68%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅8220"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.27%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'65.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.89%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.269.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.04%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.49%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.656"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.73%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -4.06%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.58%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'139.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.65%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1248"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.66%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.243"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.58%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'16.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.61%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.70%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05968"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.20%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.83%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.617"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.31%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.297"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.58%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.633"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.46%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.9833"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.53%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.779"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5899"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.91%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'FraxShare"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'6.022"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.20%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'VeChain"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.01598"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.43%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.8602"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.47%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.029.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.18%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'100.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.06%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.805.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.17%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'57.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.57%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Frax"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1.012"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.00%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'8.074"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.50%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Cronos"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.05184"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.50%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'RenderToken"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.467"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.28%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Mantle"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.4228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.32%  "
$ws.Range("E51").Style = "Normal"
